$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe out the old "diff vs. previous row" helper formulas (columns E:G,
#    plus the soon-to-exist H:J) for every data row - the new layout has no
#    formulas left in the main table.
# ---------------------------------------------------------------------------
$ws.Range("E3:J8").ClearContents()

# ---------------------------------------------------------------------------
# 2. The "Akira Red" row (old row 3) is pulled out of the main table, so the
#    remaining cursor rows shift up one row (4->3, 5->4, 6->5, 7->6).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Cursor/PointerBusy"
$ws.Range("B3").Value = 179
$ws.Range("C3").Value = 70
$ws.Range("D3").Value = 100

$ws.Range("A4").Value = "Cursor/PointerGear"
$ws.Range("B4").Value = 128
$ws.Range("C4").Value = 68
$ws.Range("D4").Value = 83

$ws.Range("A5").Value = "Cursor/PointerNewDesktop"
$ws.Range("B5").Value = 128
$ws.Range("C5").Value = 68
$ws.Range("D5").Value = 83

$ws.Range("A6").Value = "Cursor/Pointer_No_5_5"
$ws.Range("B6").Value = 357
$ws.Range("C6").Value = 82
$ws.Range("D6").Value = 69

# The old rows 7 and 8 are no longer used.
$ws.Range("A7:D8").ClearContents()

# ---------------------------------------------------------------------------
# 3. Grow Table1 from A2:G8 to A2:J7 (3 new plain "diff" columns), then give
#    the columns their new header text.
# ---------------------------------------------------------------------------
$t1 = $ws.ListObjects.Item(1)
$t1.Resize($ws.Range("A2:J7"))

$ws.Range("E2").Value = "Target Hue"
$ws.Range("F2").Value = "Target Saturation"
$ws.Range("G2").Value = "Target Value"
$ws.Range("A2").Value = "Image"
$ws.Range("H2").Value = "Hue Diff"
$ws.Range("I2").Value = "Saturation Diff"
$ws.Range("J2").Value = "Value Diff"

# ---------------------------------------------------------------------------
# 4. New Table2 holding the "Akira Red" reference color that used to live in
#    the main table.
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = "Color"
$ws.Range("M2").Value = "Hue"
$ws.Range("N2").Value = "Saturation"
$ws.Range("O2").Value = "Value"
$ws.Range("L3").Value = "Akira Red"

$t2 = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("L2:O8"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$t2.Name = "Table2"
$t2.TableStyle = "TableStyleMedium1"

# ---------------------------------------------------------------------------
# 5. Selection ends on L6, matching the author's last recorded cursor spot.
# ---------------------------------------------------------------------------
$ws.Range("L6").Select()
